# Commit: Fri, Apr 03, 2020 10:05:23 PM
#
# The table on slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") had its
# table style switched (via the PowerPoint "Table Design" style gallery)
# from the built-in "Medium Style 1 - Accent 1"
# ({6B3D3069-E90C-4F78-98FB-018C4920E2FA}) to the built-in
# "Themed Style 2 - Accent 1" ({7E80262E-1385-4CFA-9C5C-298233357D1C}).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}

$tableShape.Table.ApplyStyle("{7E80262E-1385-4CFA-9C5C-298233357D1C}")
